$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 12: Folio (column B) should be stored as a number, not text ---
$ws.Range("B12").Value = 56649874

# --- Add row 13 ---
$ws.Range("A13").Value = ""
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").Value = 140198
$ws.Range("C13").Value = "Juan"
$ws.Range("D13").Value = "Carlos"
$ws.Range("E13").Value = "Calderon"
$ws.Range("F13").Value = "Davila"
$ws.Range("G13").Value = "Director"
$ws.Range("H13").Value = "Director de la Unidad Academica"
$ws.Range("I13").Value = "23/03/2023"
$ws.Range("J13").Value = 45374
$ws.Range("J13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K13").Value = 141235
$ws.Range("L13").Value = "Niels"
$ws.Range("M13").Value = "C:/Users/MrJua/Downloads/20230318_002304.jpg"

# --- Add row 14 ---
$ws.Range("A14").Value = ""
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "14019"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "Juan"
$ws.Range("D14").Value = "Carlos"
$ws.Range("E14").Value = "Calderon"
$ws.Range("F14").Value = "Davila"
$ws.Range("G14").Value = "Director"
$ws.Range("H14").Value = "Director de la Unidad Academica"
$ws.Range("I14").Value = "23/03/2023"
$ws.Range("J14").Value = 45374
$ws.Range("J14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K14").Value = 125465
$ws.Range("L14").Value = "Niels"
$ws.Range("M14").Value = "C:/Users/MrJua/Downloads/20230318_002304.jpg"
